$d = $word.ActiveDocument

# Locate the paragraph that currently reads "10/19: " (the last paragraph
# in the access-code log before this edit).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $text = $p.Range.Text.TrimEnd([char]13)
    if ($text -eq "10/19: ") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the '10/19: ' paragraph"
}

# Append "calendar" as a brand-new paragraph first (so it lands in its own
# run rather than being merged into the "10/19: " run), then delete the
# paragraph mark that separates them so the two runs end up together in a
# single paragraph - exactly like typing "calendar" right after "10/19: ".
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$calendarPara = $d.Paragraphs($target.Index + 1)
$calendarPara.Range.InsertAfter("calendar")

$mergePoint = $d.Paragraphs($target.Index)
$mark = $d.Range($mergePoint.Range.End - 1, $mergePoint.Range.End)
$mark.Delete()

Write-Output ("merged paragraph: [" + $d.Paragraphs($target.Index).Range.Text + "]")

# Add the two new log entries ("10/24: " and "10/26: ") as their own new
# paragraphs at the end of the document.
$last = $d.Paragraphs($d.Paragraphs.Count)
$lr = $last.Range
$lr.Collapse(0)
$lr.InsertParagraphAfter()
$p1024 = $d.Paragraphs($d.Paragraphs.Count)
$p1024.Range.InsertAfter("10/24: ")

$last2 = $d.Paragraphs($d.Paragraphs.Count)
$lr2 = $last2.Range
$lr2.Collapse(0)
$lr2.InsertParagraphAfter()
$p1026 = $d.Paragraphs($d.Paragraphs.Count)
$p1026.Range.InsertAfter("10/26: ")

Write-Output ("total paragraphs now: " + $d.Paragraphs.Count)
